# Auto-generated Excel COM-interop script
# Applies updated market price / profit figures to the Leve profit tracker sheets.
# Source: scheduled runner diff against Sheets (currentAveragePrice* and derived profit columns).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 416.33334
$ws.Range("I18").Value = 416.33334
$ws.Range("K18").Value = 416.33334
$ws.Range("M18").Value = -132.33334

$ws.Range("H43").Value = 3000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H51").Value = 5500
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 8000
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = -2516
$ws.Range("N51").Value = -8968

$ws.Range("H55").Value = 251
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 2
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 2
$ws.Range("M55").Value = -286
$ws.Range("N55").Value = -430

$ws.Range("H96").Value = 208.75
$ws.Range("I96").Value = 218.28572
$ws.Range("J96").Value = 142
$ws.Range("K96").Value = 654.85716
$ws.Range("L96").Value = 426
$ws.Range("M96").Value = 718.14284
$ws.Range("N96").Value = -3172

$ws.Range("H113").Value = 6917.6665
$ws.Range("J113").Value = 8503
$ws.Range("L113").Value = 8503
$ws.Range("N113").Value = -15011

$ws.Range("H123").Value = 127593.336
$ws.Range("J123").Value = 127593.336
$ws.Range("L123").Value = 127593.336
$ws.Range("N123").Value = -137393.336

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 2725
$ws.Range("I33").Value = 2725
$ws.Range("K33").Value = 2725
$ws.Range("M33").Value = -2396

$ws.Range("H113").Value = 55555
$ws.Range("J113").Value = 55555
$ws.Range("L113").Value = 55555
$ws.Range("N113").Value = -64233

$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678

$ws.Range("H124").Value = 39000
$ws.Range("J124").Value = 39000
$ws.Range("L124").Value = 39000
$ws.Range("N124").Value = -48820

$ws.Range("H132").Value = 9193.700000000001
$ws.Range("I132").Value = 3156.1667
$ws.Range("J132").Value = 18250
$ws.Range("K132").Value = 9468.500100000001
$ws.Range("L132").Value = 54750
$ws.Range("M132").Value = -6938.500100000001
$ws.Range("N132").Value = -59810

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5100
$ws.Range("I22").Value = 5100
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 5100
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -4927
$ws.Range("N22").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 204.22223
$ws.Range("I2").Value = 261
$ws.Range("J2").Value = 5.5
$ws.Range("K2").Value = 1566
$ws.Range("L2").Value = 33
$ws.Range("M2").Value = -1453
$ws.Range("N2").Value = -259

$ws.Range("H7").Value = 22.8
$ws.Range("J7").Value = 18
$ws.Range("L7").Value = 54
$ws.Range("N7").Value = -278

$ws.Range("H8").Value = 877.5
$ws.Range("I8").Value = 877.5
$ws.Range("K8").Value = 2632.5
$ws.Range("M8").Value = -2493.5

$ws.Range("H12").Value = 112.57143
$ws.Range("J12").Value = 165.5
$ws.Range("L12").Value = 496.5
$ws.Range("N12").Value = -842.5

$ws.Range("H98").Value = 302.5
$ws.Range("I98").Value = 105
$ws.Range("K98").Value = 315
$ws.Range("M98").Value = 1183

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H126").Value = 1535.7142
$ws.Range("J126").Value = 999
$ws.Range("L126").Value = 2997
$ws.Range("N126").Value = -7937

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888

$ws.Range("H13").Value = 2746.3572
$ws.Range("I13").Value = 650
$ws.Range("K13").Value = 650
$ws.Range("M13").Value = -510

$ws.Range("H22").Value = 800
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 800
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H93").Value = 4499.8335
$ws.Range("I93").Value = 4499.8335
$ws.Range("K93").Value = 4499.8335
$ws.Range("M93").Value = -3251.8335

$ws.Range("H109").Value = 23333.334

$ws.Range("H122").Value = 4025.75
$ws.Range("J122").Value = 4066.3333
$ws.Range("L122").Value = 12198.9999
$ws.Range("N122").Value = -17098.9999

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H132").Value = 14166.5
$ws.Range("I132").Value = 8333.333000000001
$ws.Range("K132").Value = 24999.999
$ws.Range("M132").Value = -22469.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H54").Value = 36000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 36000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 36000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -37040

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H104").Value = 22370
$ws.Range("J104").Value = 22370
$ws.Range("L104").Value = 22370
$ws.Range("N104").Value = -29358

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H122").Value = 999.6667
$ws.Range("I122").Value = 999.5
$ws.Range("K122").Value = 2998.5
$ws.Range("M122").Value = -548.5

$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530
